$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 9.249492180980978
$ws.Range("C2").Value = 7.249047248251791
$ws.Range("D2").Value = 9.450550462036427
$ws.Range("E2").Value = 6.890510534567467
$ws.Range("F2").Value = 10.78533997816571
$ws.Range("G2").Value = 7.983046564903162
$ws.Range("H2").Value = 9.628621997020531
$ws.Range("I2").Value = 7.141563687632544
$ws.Range("J2").Value = 11.27794701152119
$ws.Range("K2").Value = 8.1370510003077
